$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header-row labels: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410" ---
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2404"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2410"
        }
    }
}

# --- Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U75")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
